$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'58.565.77"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.52%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.157.49"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.97%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'537.73"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.84%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'140.11"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +3.39%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.05%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.513"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +8.74%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.79%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +3.82%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.420"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +5.58%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.94%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'3.700.09"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.87%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'25.81"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +3.12%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'  +7.07%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'58.615.69"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +2.52%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = "'WrappedEther"
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = "'3.156.27"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.82%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = "'Polkadot"
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = "'6.23"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +7.01%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +5.26%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +5.97%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'373.13"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +7.56%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  +2.03%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.999"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.05%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'70.03"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +2.30%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +3.63%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +1.60%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.66%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'8.08"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +14.04%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +4.15%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value = "'PancakeSwap"
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = "'1.90"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +2.93%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = "'RenderToken"
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'6.18"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +6.25%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'21.97"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +4.84%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'5.19"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +8.68%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +5.83%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'160.10"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +0.92%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'6.23"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +4.53%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'1.37"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +13.64%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'25.28"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.59%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'2.642.40"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +9.69%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +7.19%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +4.06%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'4.18"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +4.57%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'38.80"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +5.69%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  +3.21%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  +8.79%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.999"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.04%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'3.198.02"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +2.89%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +11.61%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  +4.50%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.982"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +5.50%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'20.28"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +5.39%  "
$ws.Range('E51').Style = 'Normal'
